$d = $word.ActiveDocument

# Move to the end of the document (after the last paragraph, "first change!")
# and add a brand-new paragraph, matching the formatting already used there.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph in the document.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "after the second commit i added this line"
